$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 8334867
$ws.Range("J12").Value = 1900
$ws.Range("L12").Value = 1900
$ws.Range("N12").Value = -2240
$ws.Range("H17").Value = 3336250.2
$ws.Range("I17").Value = 1995
$ws.Range("J17").Value = 3511737.5
$ws.Range("K17").Value = 5985
$ws.Range("L17").Value = 10535212.5
$ws.Range("M17").Value = -5817
$ws.Range("N17").Value = -10535548.5
$ws.Range("H51").Value = 92594220
$ws.Range("J51").Value = 100001640
$ws.Range("L51").Value = 100001640
$ws.Range("N51").Value = -100002608
$ws.Range("H70").Value = 1675.4445
$ws.Range("J70").Value = 1869.8572
$ws.Range("L70").Value = 5609.571599999999
$ws.Range("N70").Value = -6149.571599999999
$ws.Range("H73").Value = 1675.4445
$ws.Range("J73").Value = 1869.8572
$ws.Range("L73").Value = 5609.571599999999
$ws.Range("N73").Value = -7481.571599999999
$ws.Range("H80").Value = 37037396
$ws.Range("I80").Value = 66666844
$ws.Range("J80").Value = 589.4167
$ws.Range("K80").Value = 200000532
$ws.Range("L80").Value = 1768.2501
$ws.Range("M80").Value = -199999534
$ws.Range("N80").Value = -3764.2501
$ws.Range("H83").Value = 37037396
$ws.Range("I83").Value = 66666844
$ws.Range("J83").Value = 589.4167
$ws.Range("K83").Value = 600001596
$ws.Range("L83").Value = 5304.7503
$ws.Range("M83").Value = -599996604
$ws.Range("N83").Value = -15288.7503
$ws.Range("H101").Value = 304.7143
$ws.Range("I101").Value = 226.6
$ws.Range("J101").Value = 500
$ws.Range("K101").Value = 679.8
$ws.Range("L101").Value = 1500
$ws.Range("M101").Value = 942.2
$ws.Range("N101").Value = -4744
$ws.Range("H103").Value = 661.9231
$ws.Range("J103").Value = 741.1667
$ws.Range("L103").Value = 2223.5001
$ws.Range("N103").Value = -3395.5001
$ws.Range("H106").Value = 125002984
$ws.Range("I106").Value = 142860050
$ws.Range("K106").Value = 142860050
$ws.Range("M106").Value = -142859419
$ws.Range("H107").Value = 2203.2307
$ws.Range("I107").Value = 2095.25
$ws.Range("J107").Value = 3499
$ws.Range("K107").Value = 2095.25
$ws.Range("L107").Value = 3499
$ws.Range("M107").Value = -175.25
$ws.Range("N107").Value = -7339
$ws.Range("H135").Value = 2106.0833
$ws.Range("I135").Value = 1808.1111
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 16272.9999
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -13737.9999
$ws.Range("N135").Value = -32070
$ws.Range("H137").Value = 5265011.5
$ws.Range("I137").Value = 1501
$ws.Range("J137").Value = 5884248
$ws.Range("K137").Value = 4503
$ws.Range("L137").Value = 17652744
$ws.Range("M137").Value = -1953
$ws.Range("N137").Value = -17657844
$ws.Range("H138").Value = 4968.8955
$ws.Range("J138").Value = 3136.4614
$ws.Range("L138").Value = 9409.3842
$ws.Range("N138").Value = -19689.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3769.3333
$ws.Range("I2").Value = 4012.8667
$ws.Range("K2").Value = 4012.8667
$ws.Range("M2").Value = -3899.8667
$ws.Range("H45").Value = 70699.87
$ws.Range("I45").Value = 114952.78
$ws.Range("K45").Value = 114952.78
$ws.Range("M45").Value = -114575.78
$ws.Range("H61").Value = 1734975.4
$ws.Range("I61").Value = 4608.1035
$ws.Range("J61").Value = 11771106
$ws.Range("K61").Value = 4608.1035
$ws.Range("L61").Value = 11771106
$ws.Range("M61").Value = -4396.1035
$ws.Range("N61").Value = -11771530
$ws.Range("H74").Value = 1013561.2
$ws.Range("I74").Value = 1719.5834
$ws.Range("K74").Value = 1719.5834
$ws.Range("M74").Value = -845.5834
$ws.Range("H77").Value = 1013561.2
$ws.Range("I77").Value = 1719.5834
$ws.Range("K77").Value = 8597.916999999999
$ws.Range("M77").Value = -4229.916999999999
$ws.Range("H116").Value = 3769.3333
$ws.Range("I116").Value = 4012.8667
$ws.Range("K116").Value = 4012.8667
$ws.Range("M116").Value = -1718.8667
$ws.Range("H136").Value = 1734975.4
$ws.Range("I136").Value = 4608.1035
$ws.Range("J136").Value = 11771106
$ws.Range("K136").Value = 13824.3105
$ws.Range("L136").Value = 35313318
$ws.Range("M136").Value = -11274.3105
$ws.Range("N136").Value = -35318418

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3769.3333
$ws.Range("I3").Value = 4012.8667
$ws.Range("K3").Value = 4012.8667
$ws.Range("M3").Value = -3898.8667
$ws.Range("H107").Value = 12076.5
$ws.Range("I107").Value = 13106.154
$ws.Range("K107").Value = 13106.154
$ws.Range("M107").Value = -11186.154
$ws.Range("I134").Value = 3054.077
$ws.Range("K134").Value = 9162.231
$ws.Range("M134").Value = -6627.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3317.66
$ws.Range("J31").Value = 3283.3264
$ws.Range("L31").Value = 3283.3264
$ws.Range("N31").Value = -3873.3264
$ws.Range("H34").Value = 3317.66
$ws.Range("J34").Value = 3283.3264
$ws.Range("L34").Value = 3283.3264
$ws.Range("N34").Value = -3687.3264
$ws.Range("H86").Value = 35694.3
$ws.Range("J86").Value = 46283
$ws.Range("L86").Value = 46283
$ws.Range("N86").Value = -48529
$ws.Range("H89").Value = 35694.3
$ws.Range("J89").Value = 46283
$ws.Range("L89").Value = 231415
$ws.Range("N89").Value = -242647
$ws.Range("H99").Value = 66206.75
$ws.Range("I99").Value = 4913
$ws.Range("K99").Value = 4913
$ws.Range("M99").Value = -3415
$ws.Range("H126").Value = 66206.75
$ws.Range("I126").Value = 4913
$ws.Range("K126").Value = 14739
$ws.Range("M126").Value = -12269
$ws.Range("H132").Value = 18521790
$ws.Range("I132").Value = 3371.3333
$ws.Range("K132").Value = 10113.9999
$ws.Range("M132").Value = -7583.999899999999
$ws.Range("H141").Value = 643332.8
$ws.Range("J141").Value = 643332.8
$ws.Range("L141").Value = 643332.8
$ws.Range("N141").Value = -653692.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 447
$ws.Range("I2").Value = 93.5
$ws.Range("J2").Value = 623.75
$ws.Range("K2").Value = 561
$ws.Range("L2").Value = 3742.5
$ws.Range("M2").Value = -448
$ws.Range("N2").Value = -3968.5
$ws.Range("H95").Value = 24995
$ws.Range("J95").Value = 24995
$ws.Range("L95").Value = 74985
$ws.Range("N95").Value = -79103
$ws.Range("H97").Value = 500
$ws.Range("I97").Value = 500
$ws.Range("K97").Value = 1500
$ws.Range("M97").Value = -1004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I132").Value = 1556.238
$ws.Range("K132").Value = 4668.714
$ws.Range("M132").Value = -2138.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10370.177
$ws.Range("J7").Value = 17105.555
$ws.Range("L7").Value = 17105.555
$ws.Range("N7").Value = -17329.555
$ws.Range("H40").Value = 4623.533
$ws.Range("I40").Value = 3993
$ws.Range("J40").Value = 5043.8887
$ws.Range("K40").Value = 3993
$ws.Range("L40").Value = 5043.8887
$ws.Range("M40").Value = -3857
$ws.Range("N40").Value = -5315.8887
$ws.Range("H68").Value = 3085.75
$ws.Range("I68").Value = 3181.3333
$ws.Range("J68").Value = 2799
$ws.Range("K68").Value = 3181.3333
$ws.Range("L68").Value = 2799
$ws.Range("M68").Value = -2432.3333
$ws.Range("N68").Value = -4297
$ws.Range("H71").Value = 3085.75
$ws.Range("I71").Value = 3181.3333
$ws.Range("J71").Value = 2799
$ws.Range("K71").Value = 15906.6665
$ws.Range("L71").Value = 13995
$ws.Range("M71").Value = -12162.6665
$ws.Range("N71").Value = -21483
$ws.Range("H126").Value = 10370.177
$ws.Range("J126").Value = 17105.555
$ws.Range("L126").Value = 51316.665
$ws.Range("N126").Value = -56256.665
$ws.Range("H132").Value = 4893.25
$ws.Range("I132").Value = 4066.125
$ws.Range("J132").Value = 5996.0835
$ws.Range("K132").Value = 12198.375
$ws.Range("L132").Value = 17988.2505
$ws.Range("M132").Value = -9668.375
$ws.Range("N132").Value = -23048.2505
$ws.Range("H136").Value = 3421.158
$ws.Range("I136").Value = 1650.2
$ws.Range("J136").Value = 5388.8887
$ws.Range("K136").Value = 4950.6
$ws.Range("L136").Value = 16166.6661
$ws.Range("M136").Value = -2400.6
$ws.Range("N136").Value = -21266.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3240708.5
$ws.Range("I81").Value = 4124089
$ws.Range("K81").Value = 8248178
$ws.Range("M81").Value = -8247117
$ws.Range("H84").Value = 3240708.5
$ws.Range("I84").Value = 4124089
$ws.Range("K84").Value = 41240890
$ws.Range("M84").Value = -41235586
$ws.Range("H122").Value = 3042.818
$ws.Range("I122").Value = 3042.818
$ws.Range("K122").Value = 9128.454000000002
$ws.Range("M122").Value = -6678.454000000002
$ws.Range("H132").Value = 38770.668
$ws.Range("I132").Value = 49207.332
$ws.Range("K132").Value = 147621.996
$ws.Range("M132").Value = -145091.996
$ws.Range("H136").Value = 29155.611
$ws.Range("I136").Value = 37950.926
$ws.Range("K136").Value = 113852.778
$ws.Range("M136").Value = -111302.778
